$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (the existing D:K data columns shift to E:L),
# matching the updated report that adds one more period of financial data.
$ws.Columns.Item(4).Insert()

# Copy number formats/styles from column E (the former column D) into the
# newly inserted column D so the new cells inherit the same look (date
# format on row headers, number format on data rows, etc.) and reuse the
# existing style entries instead of creating new ones.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly inserted column D with the new period's values.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 3666800
$ws.Range("D9").Value = 2534000
$ws.Range("D10").Value = 1132800
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 77600
$ws.Range("D15").Value = 88100
$ws.Range("D17").Value = 3429800
$ws.Range("D18").Value = 237000
$ws.Range("D20").Value = -10100
$ws.Range("D21").Value = 368700
$ws.Range("D22").Value = 44100
$ws.Range("D23").Value = 182800
$ws.Range("D24").Value = 11500
$ws.Range("D26").Value = 171400
$ws.Range("D27").Value = 157100
$ws.Range("D29").Value = -16900
$ws.Range("D32").Value = 10100
$ws.Range("D33").Value = 140200
$ws.Range("D35").Value = 140200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 245000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 989400
$ws.Range("D44").Value = 496500
$ws.Range("D45").Value = 227500
$ws.Range("D46").Value = 1958400
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 503300
$ws.Range("D49").Value = 3589500
$ws.Range("D52").Value = 552600
$ws.Range("D54").Value = 6603900
$ws.Range("D57").Value = 640700
$ws.Range("D58").Value = 6300
$ws.Range("D59").Value = 552300
$ws.Range("D60").Value = 1199300
$ws.Range("D61").Value = 1192400
$ws.Range("D62").Value = 735200
$ws.Range("D66").Value = 3334100
$ws.Range("D72").Value = 991800
$ws.Range("D76").Value = 3269800
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 140200
$ws.Range("D83").Value = 141900
$ws.Range("D89").Value = 226400
$ws.Range("D91").Value = -69600
$ws.Range("D94").Value = -167900
$ws.Range("D96").Value = 0
$ws.Range("D100").Value = -47200
$ws.Range("D101").Value = -28400
$ws.Range("D102").Value = -17000

$wb.Save()
